$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "body" row (row 23) ---------------------------------------
# Start from A11's current formatting (it already reuses the fillId=8 "light"
# cellXf) so the new row shares the same fill, then tweak the font to an
# italic, 12pt body-copy style.
$ws.Range("A11").Copy($ws.Range("A23"))
$ws.Range("A23").Value = "body"
$ws.Range("B23").Value = 11

$bodyFont = $ws.Range("A23").Font
$bodyFont.Size = 12
$bodyFont.Italic = $true

# --- Fix A11's font colour: was reading as a raw theme number --------------
# xlThemeColorLight2 (4) resolves to the workbook's theme="2" colour.
$ws.Range("A11").Font.ThemeColor = 4

# --- Restore the selection back to the row that was just fixed -------------
$ws.Range("A11").Select()
